# Player.xlsx - "Property" sheet (sheet1 / Worksheets.Item(1))
#
# 1) Move the active selection to H78 (was E76).
# 2) Flip the "Save" column (E) from TRUE to FALSE for rows 44-67.
# 3) Rows 76-77 were highlighted (red text on yellow fill, plus a text
#    number format) as "new" rows; the edit removes that highlighting,
#    restoring plain/default formatting. Row 76 col A keeps a Text
#    number format (no color/fill), everything else goes back to the
#    workbook's default style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1) Selection -----------------------------------------------------
$ws.Range("H78").Select()

# --- 2) Save column TRUE -> FALSE for rows 44..67 ----------------------
for ($row = 44; $row -le 67; $row++) {
    $ws.Range("E$row").Value = $false
}

# --- 3) Clear the "new row" highlight formatting -----------------------
# Row 76: column A keeps the Text number format, but loses the red
# font / yellow fill; the rest of the row drops back to plain default
# formatting.
$ws.Range("A76").ClearFormats()
$ws.Range("A76").NumberFormat = "@"
$ws.Range("B76").ClearFormats()
$ws.Range("G76").ClearFormats()
$ws.Range("H76").ClearFormats()
$ws.Range("I76").ClearFormats()
$ws.Range("J76").ClearFormats()

# Row 77: entire row drops back to plain default formatting.
$ws.Range("A77").ClearFormats()
$ws.Range("B77").ClearFormats()
$ws.Range("G77").ClearFormats()
$ws.Range("H77").ClearFormats()
$ws.Range("I77").ClearFormats()
$ws.Range("J77").ClearFormats()
